$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7459.8623
$ws.Range("I19").Value = 1702.9412
$ws.Range("J19").Value = 15615.5
$ws.Range("K19").Value = 1702.9412
$ws.Range("L19").Value = 15615.5
$ws.Range("M19").Value = -1527.9412
$ws.Range("N19").Value = -15965.5
$ws.Range("H62").Value = 4824.636
$ws.Range("I62").Value = 4674.5557
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 4674.5557
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -4050.5557
$ws.Range("N62").Value = -6748
$ws.Range("H65").Value = 4824.636
$ws.Range("I65").Value = 4674.5557
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 23372.7785
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -20252.7785
$ws.Range("N65").Value = -33740
$ws.Range("H68").Value = 59999
$ws.Range("J68").Value = 59999
$ws.Range("L68").Value = 59999
$ws.Range("N68").Value = -61497
$ws.Range("H71").Value = 59999
$ws.Range("J71").Value = 59999
$ws.Range("L71").Value = 179997
$ws.Range("N71").Value = -187485
$ws.Range("H112").Value = 2358.76
$ws.Range("J112").Value = 2465.3333
$ws.Range("L112").Value = 7395.999899999999
$ws.Range("N112").Value = -9611.999899999999
$ws.Range("H137").Value = 4080.7258
$ws.Range("I137").Value = 4000.5
$ws.Range("J137").Value = 4829.5
$ws.Range("K137").Value = 12001.5
$ws.Range("L137").Value = 14488.5
$ws.Range("M137").Value = -9451.5
$ws.Range("N137").Value = -19588.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7302.9297
$ws.Range("I32").Value = 4770.216
$ws.Range("K32").Value = 4770.216
$ws.Range("M32").Value = -4483.216
$ws.Range("H61").Value = 3751.4614
$ws.Range("J61").Value = 6497
$ws.Range("L61").Value = 6497
$ws.Range("N61").Value = -6921
$ws.Range("H132").Value = 804.48
$ws.Range("I132").Value = 688.0417
$ws.Range("K132").Value = 2064.1251
$ws.Range("M132").Value = 465.8748999999998
$ws.Range("H136").Value = 3751.4614
$ws.Range("J136").Value = 6497
$ws.Range("L136").Value = 19491
$ws.Range("N136").Value = -24591

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1888.3889
$ws.Range("I86").Value = 1056
$ws.Range("J86").Value = 2720.7778
$ws.Range("K86").Value = 1056
$ws.Range("L86").Value = 2720.7778
$ws.Range("M86").Value = 67
$ws.Range("N86").Value = -4966.7778
$ws.Range("H89").Value = 1888.3889
$ws.Range("I89").Value = 1056
$ws.Range("J89").Value = 2720.7778
$ws.Range("K89").Value = 5280
$ws.Range("L89").Value = 13603.889
$ws.Range("M89").Value = 336
$ws.Range("N89").Value = -24835.889
$ws.Range("H134").Value = 5142.976
$ws.Range("I134").Value = 4855.758
$ws.Range("J134").Value = 6196.1113
$ws.Range("K134").Value = 14567.274
$ws.Range("L134").Value = 18588.3339
$ws.Range("M134").Value = -12032.274
$ws.Range("N134").Value = -23658.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2911.8928
$ws.Range("I31").Value = 2077.5
$ws.Range("J31").Value = 5971.3335
$ws.Range("K31").Value = 2077.5
$ws.Range("L31").Value = 5971.3335
$ws.Range("M31").Value = -1782.5
$ws.Range("N31").Value = -6561.3335
$ws.Range("H34").Value = 2911.8928
$ws.Range("I34").Value = 2077.5
$ws.Range("J34").Value = 5971.3335
$ws.Range("K34").Value = 2077.5
$ws.Range("L34").Value = 5971.3335
$ws.Range("M34").Value = -1875.5
$ws.Range("N34").Value = -6375.3335
$ws.Range("H132").Value = 9762
$ws.Range("I132").Value = 7741.091
$ws.Range("K132").Value = 23223.273
$ws.Range("M132").Value = -20693.273
$ws.Range("H134").Value = 2154.8071
$ws.Range("I134").Value = 1960.1666
$ws.Range("K134").Value = 5880.4998
$ws.Range("M134").Value = -3345.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 53099.21
$ws.Range("I4").Value = 282.08334
$ws.Range("J4").Value = 143642.86
$ws.Range("K4").Value = 846.2500200000001
$ws.Range("L4").Value = 430928.58
$ws.Range("M4").Value = -734.2500200000001
$ws.Range("N4").Value = -431152.58
$ws.Range("H37").Value = 110142.29
$ws.Range("J37").Value = 110142.29
$ws.Range("L37").Value = 330426.87
$ws.Range("N37").Value = -330650.87
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H97").Value = 2733.5
$ws.Range("I97").Value = 2592.889
$ws.Range("J97").Value = 3999
$ws.Range("K97").Value = 7778.667
$ws.Range("L97").Value = 11997
$ws.Range("M97").Value = -7282.667
$ws.Range("N97").Value = -12989
$ws.Range("H114").Value = 740.625
$ws.Range("I114").Value = 759.4
$ws.Range("K114").Value = 2278.2
$ws.Range("M114").Value = 975.8000000000002
$ws.Range("H117").Value = 1281.5
$ws.Range("J117").Value = 2532
$ws.Range("L117").Value = 7596
$ws.Range("N117").Value = -14480
$ws.Range("H129").Value = 1194
$ws.Range("I129").Value = 741.9091
$ws.Range("J129").Value = 2437.25
$ws.Range("K129").Value = 2225.7273
$ws.Range("L129").Value = 7311.75
$ws.Range("M129").Value = 2774.2727
$ws.Range("N129").Value = -17311.75
$ws.Range("H131").Value = 489538.3
$ws.Range("J131").Value = 2606.9216
$ws.Range("L131").Value = 7820.764800000001
$ws.Range("N131").Value = -17900.7648

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 25266498
$ws.Range("I80").Value = 36925516
$ws.Range("J80").Value = 5296.5
$ws.Range("K80").Value = 36925516
$ws.Range("L80").Value = 5296.5
$ws.Range("M80").Value = -36924518
$ws.Range("N80").Value = -7292.5
$ws.Range("H83").Value = 25266498
$ws.Range("I83").Value = 36925516
$ws.Range("J83").Value = 5296.5
$ws.Range("K83").Value = 184627580
$ws.Range("L83").Value = 26482.5
$ws.Range("M83").Value = -184622588
$ws.Range("N83").Value = -36466.5
$ws.Range("H126").Value = 8832.4
$ws.Range("I126").Value = 9299.333000000001
$ws.Range("J126").Value = 8632.286
$ws.Range("K126").Value = 27897.999
$ws.Range("L126").Value = 25896.858
$ws.Range("M126").Value = -25427.999
$ws.Range("N126").Value = -30836.858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3237.7856
$ws.Range("I40").Value = 2792.111
$ws.Range("J40").Value = 4040
$ws.Range("K40").Value = 2792.111
$ws.Range("L40").Value = 4040
$ws.Range("M40").Value = -2656.111
$ws.Range("N40").Value = -4312
$ws.Range("H122").Value = 3410.4167
$ws.Range("I122").Value = 2787.6
$ws.Range("K122").Value = 8362.799999999999
$ws.Range("M122").Value = -5912.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 120465
$ws.Range("J141").Value = 120465
$ws.Range("L141").Value = 120465
$ws.Range("N141").Value = -130825
